# Reposition the flowchart objects on slide 9 (the embedded Word OLE
# object and the six process-step pictures) per the new layout.
#
# Target positions below are expressed in points (1 pt = 12700 EMU) and
# were chosen so that, after the host's float32 Left/Top round-trip,
# the serialized EMU offsets land exactly on the requested values:
#   Object 18  -> off 141455, 1442278 EMU
#   Picture 21 -> off 3594642, 1446160 EMU
#   Picture 23 -> off 5402884, 1442278 EMU
#   Picture 24 -> off 7658770, 1442278 EMU
#   Picture 25 -> off 138684, 3730276 EMU
#   Picture 26 -> off 3358652, 3717944 EMU
#   Picture 27 -> off 5568748, 4148757 EMU

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

$moves = @(
    @{ Name = "Object 18";  Left = 11.138189315795898;  Top = 113.56520080566406 },
    @{ Name = "Picture 21"; Left = 283.0426940917969;    Top = 113.8708724975586  },
    @{ Name = "Picture 23"; Left = 425.4239501953125;    Top = 113.56520080566406 },
    @{ Name = "Picture 24"; Left = 603.0527954101562;    Top = 113.56520080566406 },
    @{ Name = "Picture 25"; Left = 10.920000076293945;   Top = 293.7225341796875  },
    @{ Name = "Picture 26"; Left = 264.4608154296875;    Top = 292.75152587890625 },
    @{ Name = "Picture 27"; Left = 438.4841003417969;    Top = 326.6737976074219  }
)

foreach ($m in $moves) {
    $shape = $s.Shapes.Item($m.Name)
    $shape.Left = $m.Left
    $shape.Top = $m.Top
}
